$wb = $excel.ActiveWorkbook

# --- FatosIn sheet: fill in the new "tipo" column (N) / value column (O) ---
# for the 7 newly-added element rows (36-42).
$ws = $wb.Worksheets.Item("FatosIn")

$tipos = @{
    36 = '"Laje_Básica"'
    37 = '"Coluna_20x20"'
    38 = '"Barra_b1"'
    39 = '"Barra_P1"'
    40 = '"Barra_N1"'
    41 = '"Barra_E1"'
    42 = '"Acop_P1"'
}

foreach ($r in 36..42) {
    $ws.Range("N$r").Value2 = "tipo"
    $ws.Range("O$r").Value2 = $tipos[$r]
}

# --- Re-align column O (the new value column) to left instead of center,
# matching the rest of the "value" columns' look now that it actually holds text.
$ws.Range("O1").HorizontalAlignment = -4131
$ws.Range("O2:O42").HorizontalAlignment = -4131

# --- View state: the workbook was left with FatosIn as the active/visible
# sheet, with row 36 selected.
$ws.Activate()
$ws.Rows.Item(36).Select()
